# Apply updated dSF (column F) values, as described in the commit:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 1
    4  = 3
    5  = 4
    6  = 0
    7  = -3
    8  = -3
    9  = -2
    10 = -1
    11 = -2
    12 = -1
    14 = -3
    15 = -2
    16 = -1
    17 = -3
    18 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
